$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.Text = "Quel est le plus long fleuve du monde ?"
$d.Paragraphs.Item(2).Range.Text = "Le Nil"
$d.Paragraphs.Item(3).Range.Text = "L'Amazone"
$d.Paragraphs.Item(4).Range.Text = "Le Yang-Tsé-Kiang"
$d.Paragraphs.Item(5).Range.Text = "Le Mississippi"
$d.Paragraphs.Item(7).Range.Text = "Quel est le plus haut sommet du monde ?"
$d.Paragraphs.Item(8).Range.Text = "L'Everest"
$d.Paragraphs.Item(9).Range.Text = "Le K2"
$d.Paragraphs.Item(10).Range.Text = "Le Mont Blanc"
$d.Paragraphs.Item(11).Range.Text = "Le Kilimandjaro"
$d.Paragraphs.Item(13).Range.Text = "Quelle est la capitale de l'Argentine ?"
$d.Paragraphs.Item(14).Range.Text = "Rio de Janeiro"
$d.Paragraphs.Item(15).Range.Text = "Buenos Aires"
$d.Paragraphs.Item(16).Range.Text = "Santiago"
$d.Paragraphs.Item(17).Range.Text = "Lima"
$d.Paragraphs.Item(19).Range.Text = "Quel est le plus grand récif corallien du monde ?"
$d.Paragraphs.Item(20).Range.Text = "La Grande Barrière de corail"
$d.Paragraphs.Item(21).Range.Text = "Le récif de la mer Rouge"
$d.Paragraphs.Item(22).Range.Text = "Le récif de la Nouvelle-Calédonie"
$d.Paragraphs.Item(23).Range.Text = "Le récif de la Mésoamérique"
$d.Paragraphs.Item(25).Range.Text = "Quelle est la capitale du Japon ?"
$d.Paragraphs.Item(26).Range.Text = "Kyoto"
$d.Paragraphs.Item(27).Range.Text = "Osaka"
$d.Paragraphs.Item(28).Range.Text = "Tokyo"
$d.Paragraphs.Item(29).Range.Text = "Hiroshima"
$d.Paragraphs.Item(31).Range.Text = "Quel est le nom du plus célèbre détective créé par Arthur Conan Doyle ?"
$d.Paragraphs.Item(32).Range.Text = "Sherlock Holmes"
$d.Paragraphs.Item(33).Range.Text = "Hercule Poirot"
$d.Paragraphs.Item(34).Range.Text = "Miss Marple"
$d.Paragraphs.Item(35).Range.Text = "Arsène Lupin"
$d.Paragraphs.Item(37).Range.Text = "Quel est le plus grand arbre du monde ?"
$d.Paragraphs.Item(38).Range.Text = "Le séquoia géant"
$d.Paragraphs.Item(39).Range.Text = "Le séquoia à feuilles d'if"
$d.Paragraphs.Item(40).Range.Text = "Le baobab"
$d.Paragraphs.Item(41).Range.Text = "Le chêne vert"
$d.Paragraphs.Item(43).Range.Text = "Quel est le plus grand insecte du monde ?"
$d.Paragraphs.Item(44).Range.Text = "Le scarabée Goliath"
$d.Paragraphs.Item(45).Range.Text = "Le papillon atlas"
$d.Paragraphs.Item(46).Range.Text = "Le phasme géant"
$d.Paragraphs.Item(47).Range.Text = "La fourmi géante"
$d.Paragraphs.Item(49).Range.Text = "Quel est le plus grand volcan du monde ?"
$d.Paragraphs.Item(50).Range.Text = "Le mont Vésuve"
$d.Paragraphs.Item(51).Range.Text = "Le mont Fuji"
$d.Paragraphs.Item(52).Range.Text = "Le Mauna Loa"
$d.Paragraphs.Item(53).Range.Text = "Le mont Saint Helens"
$d.Paragraphs.Item(55).Range.Text = "Quel est le nom du premier homme à avoir marché sur la Lune ?"
$d.Paragraphs.Item(56).Range.Text = "Neil Armstrong"
$d.Paragraphs.Item(57).Range.Text = "Buzz Aldrin"
$d.Paragraphs.Item(58).Range.Text = "Youri Gagarine"
$d.Paragraphs.Item(59).Range.Text = "Alan Shepard"
$d.Paragraphs.Item(61).Range.Text = "Quel est le plus grand pays du monde en superficie ?"
$d.Paragraphs.Item(62).Range.Text = "La Russie"
$d.Paragraphs.Item(63).Range.Text = "Le Canada"
$d.Paragraphs.Item(64).Range.Text = "La Chine"
$d.Paragraphs.Item(65).Range.Text = "Les États-Unis"
$d.Paragraphs.Item(67).Range.Text = "﻿Quel est le plus grand océan du monde ?"
$d.Paragraphs.Item(68).Range.Text = "L'océan Atlantique"
$d.Paragraphs.Item(69).Range.Text = "L'océan Indien"
$d.Paragraphs.Item(70).Range.Text = "L'océan Pacifique"
$d.Paragraphs.Item(71).Range.Text = "L'océan Arctique"
$d.Paragraphs.Item(73).Range.Text = "Quel est le plus grand désert du monde ?"
$d.Paragraphs.Item(74).Range.Text = "Le Sahara"
$d.Paragraphs.Item(75).Range.Text = "Le désert de Gobi"
$d.Paragraphs.Item(76).Range.Text = "Le désert d'Atacama"
$d.Paragraphs.Item(77).Range.Text = "Le désert Antarctique"
$d.Paragraphs.Item(79).Range.Text = "Quel est le plus grand lac du monde ?"
$d.Paragraphs.Item(80).Range.Text = "La mer Caspienne"
$d.Paragraphs.Item(81).Range.Text = "Le lac Supérieur"
$d.Paragraphs.Item(82).Range.Text = "Le lac Victoria"
$d.Paragraphs.Item(83).Range.Text = "Le lac Baïkal"
$d.Paragraphs.Item(85).Range.Text = "Qui a peint la Joconde ?"
$d.Paragraphs.Item(86).Range.Text = "Vincent van Gogh"
$d.Paragraphs.Item(87).Range.Text = "Léonard de Vinci"
$d.Paragraphs.Item(88).Range.Text = "Michel-Ange"
$d.Paragraphs.Item(89).Range.Text = "Raphaël"
$d.Paragraphs.Item(91).Range.Text = "Quelle est la capitale de l'Australie ?"
$d.Paragraphs.Item(92).Range.Text = "Sydney"
$d.Paragraphs.Item(93).Range.Text = "Melbourne"
$d.Paragraphs.Item(94).Range.Text = "Canberra"
$d.Paragraphs.Item(95).Range.Text = "Perth"
$d.Paragraphs.Item(97).Range.Text = "Quel est le plus grand amphibien du monde ?"
$d.Paragraphs.Item(98).Range.Text = "La salamandre géante de Chine"
$d.Paragraphs.Item(99).Range.Text = "La salamandre géante du Japon"
$d.Paragraphs.Item(100).Range.Text = "La grenouille taureau"
$d.Paragraphs.Item(101).Range.Text = "Le triton"
$d.Paragraphs.Item(103).Range.Text = "Quel est le plus grand animal terrestre ?"
$d.Paragraphs.Item(104).Range.Text = "La baleine bleue"
$d.Paragraphs.Item(105).Range.Text = "L'éléphant d'Afrique"
$d.Paragraphs.Item(106).Range.Text = "Le rhinocéros"
$d.Paragraphs.Item(107).Range.Text = "L'hippopotame"
$d.Paragraphs.Item(109).Range.Text = "Qui a écrit `"Les Misérables`" ?"
$d.Paragraphs.Item(110).Range.Text = "Gustave Flaubert"
$d.Paragraphs.Item(111).Range.Text = "Victor Hugo"
$d.Paragraphs.Item(112).Range.Text = "Émile Zola"
$d.Paragraphs.Item(113).Range.Text = "Alexandre Dumas"
$d.Paragraphs.Item(115).Range.Text = "Quel est le plus grand poisson du monde ?"
$d.Paragraphs.Item(116).Range.Text = "Le requin blanc"
$d.Paragraphs.Item(117).Range.Text = "Le requin-baleine"
$d.Paragraphs.Item(118).Range.Text = "Le poisson-lune"
$d.Paragraphs.Item(119).Range.Text = "Le thon rouge"
